$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Daño foliar -> Leaf_dmg)
$ws.Name = "Leaf_dmg"

# Translate / rename the header row (also updates the structured table's column names)
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Field"
$ws.Range("C1").Value = "Treatment"
$ws.Range("D1").Value = "Repeat"
$ws.Range("E1").Value = "Leaves_dmg_10leaves"
$ws.Range("F1").Value = "Marks_5leaves"
$ws.Range("G1").Value = "Observations"

# Resize columns B:F (drops their bestFit auto-size in favour of explicit widths)
$ws.Columns.Item(2).ColumnWidth = 7
$ws.Columns.Item(3).ColumnWidth = 12
$ws.Columns.Item(4).ColumnWidth = 9
$ws.Columns.Item(5).ColumnWidth = 22.6666666666667
$ws.Columns.Item(6).ColumnWidth = 17.3333333333333

# Move the active selection to G11
$ws.Range("G11").Select() | Out-Null
